$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.75"
$ws.Range("E2").Value = "'-4.71%"
$ws.Range("D3").Value = "'49.00"
$ws.Range("E3").Value = "'-1.15%"
$ws.Range("D4").Value = "'5.175"
$ws.Range("E4").Value = "'-2.73%"
$ws.Range("D5").Value = "'0.07714"
$ws.Range("E5").Value = "'-5.43%"
$ws.Range("D6").Value = "'4.511"
$ws.Range("E6").Value = "'-2.05%"
$ws.Range("D7").Value = "'1.340"
$ws.Range("E7").Value = "'13.62%"
$ws.Range("D8").Value = "'1.545"
$ws.Range("E8").Value = "'-7.55%"
$ws.Range("D9").Value = "'0.1231"
$ws.Range("E9").Value = "'-9.15%"
$ws.Range("D10").Value = "'0.1931"
$ws.Range("E10").Value = "'-1.38%"
$ws.Range("D11").Value = "'0.04672"
$ws.Range("E11").Value = "'3.29%"
$ws.Range("D12").Value = "'0.09399"
$ws.Range("E12").Value = "'-2.57%"
$ws.Range("E13").Value = "'0.01%"
$ws.Range("D14").Value = "'0.001265"
$ws.Range("E14").Value = "'-4.45%"
$ws.Range("D15").Value = "'0.04175"
$ws.Range("E15").Value = "'-3.07%"
$ws.Range("D16").Value = "'0.005828"
$ws.Range("E16").Value = "'-2.23%"
$ws.Range("E17").Value = "'-2.07%"
$ws.Range("E18").Value = "'-6.76%"
$ws.Range("E19").Value = "'2.73%"
$ws.Range("D20").Value = "'7.958"
$ws.Range("E20").Value = "'-2.15%"
$ws.Range("D21").Value = "'0.1339"
$ws.Range("E21").Value = "'-5.69%"
$ws.Range("E22").Value = "'-0.42%"
$ws.Range("E23").Value = "'-2.31%"
$ws.Range("D24").Value = "'0.004051"
$ws.Range("E24").Value = "'-5.07%"
$ws.Range("D25").Value = "'0.0001353"
$ws.Range("E25").Value = "'0.24%"
$ws.Range("E26").Value = "'0.93%"
$ws.Range("D38").Value = "'0.02576"
$ws.Range("E38").Value = "'-7.31%"
$ws.Range("D39").Value = "'0.05795"
$ws.Range("E39").Value = "'3.71%"
$ws.Range("E40").Value = "'70.78%"
$ws.Range("D41").Value = "'0.007963"
$ws.Range("E41").Value = "'3.11%"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("E42").Value = "'-1.88%"
$ws.Range("D43").Value = "'0.008377"
$ws.Range("E43").Value = "'8.93%"
$ws.Range("D44").Value = "'0.007674"
$ws.Range("E44").Value = "'-5.08%"
$ws.Range("D45").Value = "'0.3370"
$ws.Range("E45").Value = "'-4.17%"
$ws.Range("D46").Value = "'0.00007016"
$ws.Range("E46").Value = "'2.82%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.25%"
$ws.Range("E48").Value = "'-7.55%"
$ws.Range("E49").Value = "'0.17%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.25%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.25%"
